# Rename the "_old" / "_new" column-header suffixes to the respective
# input-file format versions ("_FV2410" / "_FV2504"), wrap the header/data
# range in an Excel Table (ListObject), and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells in row 1 (A1:U1) -----------------------
# Columns A-J + L-U used to carry a "_old"/"_new" suffix; K ("diff") is
# left untouched.
$headerMap = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# --- 2. Turn the used range into an Excel Table (ListObject) -----------
$dataRange = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (row 1) -----------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
